$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 258.68
$ws.Cells.Item(15, 9).Value = 258.68
$ws.Cells.Item(15, 11).Value = 776.04
$ws.Cells.Item(15, 13).Value = -607.04
$ws.Cells.Item(111, 8).Value = 14341.667
$ws.Cells.Item(111, 9).Value = 3634.375
$ws.Cells.Item(111, 10).Value = 100000
$ws.Cells.Item(111, 11).Value = 10903.125
$ws.Cells.Item(111, 12).Value = 300000
$ws.Cells.Item(111, 13).Value = -7836.125
$ws.Cells.Item(111, 14).Value = -306134
$ws.Cells.Item(132, 8).Value = 1975138.6
$ws.Cells.Item(132, 9).Value = 399245.2
$ws.Cells.Item(132, 11).Value = 1197735.6
$ws.Cells.Item(132, 13).Value = -1195205.6
$ws.Cells.Item(137, 8).Value = 17535066
$ws.Cells.Item(137, 9).Value = 4465193
$ws.Cells.Item(137, 10).Value = 50803836
$ws.Cells.Item(137, 11).Value = 13395579
$ws.Cells.Item(137, 12).Value = 152411508
$ws.Cells.Item(137, 13).Value = -13393029
$ws.Cells.Item(137, 14).Value = -152416608

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10635.35
$ws.Cells.Item(32, 9).Value = 8062.898
$ws.Cells.Item(32, 10).Value = 29500
$ws.Cells.Item(32, 11).Value = 8062.898
$ws.Cells.Item(32, 12).Value = 29500
$ws.Cells.Item(32, 13).Value = -7775.898
$ws.Cells.Item(32, 14).Value = -30074
$ws.Cells.Item(62, 8).Value = 7500
$ws.Cells.Item(62, 10).Value = 7500
$ws.Cells.Item(62, 12).Value = 7500
$ws.Cells.Item(62, 14).Value = -8748
$ws.Cells.Item(65, 8).Value = 7500
$ws.Cells.Item(65, 10).Value = 7500
$ws.Cells.Item(65, 12).Value = 22500
$ws.Cells.Item(65, 14).Value = -28740
$ws.Cells.Item(132, 8).Value = 27973482
$ws.Cells.Item(132, 9).Value = 30055474
$ws.Cells.Item(132, 10).Value = 17860950
$ws.Cells.Item(132, 11).Value = 90166422
$ws.Cells.Item(132, 12).Value = 53582850
$ws.Cells.Item(132, 13).Value = -90163892
$ws.Cells.Item(132, 14).Value = -53587910

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1955.52
$ws.Cells.Item(86, 9).Value = 1956.5858
$ws.Cells.Item(86, 10).Value = 1850
$ws.Cells.Item(86, 11).Value = 1956.5858
$ws.Cells.Item(86, 12).Value = 1850
$ws.Cells.Item(86, 13).Value = -833.5858000000001
$ws.Cells.Item(86, 14).Value = -4096
$ws.Cells.Item(89, 8).Value = 1955.52
$ws.Cells.Item(89, 9).Value = 1956.5858
$ws.Cells.Item(89, 10).Value = 1850
$ws.Cells.Item(89, 11).Value = 9782.929
$ws.Cells.Item(89, 12).Value = 9250
$ws.Cells.Item(89, 13).Value = -4166.929
$ws.Cells.Item(89, 14).Value = -20482
$ws.Cells.Item(134, 8).Value = 13318952
$ws.Cells.Item(134, 9).Value = 15262357
$ws.Cells.Item(134, 11).Value = 45787071
$ws.Cells.Item(134, 13).Value = -45784536
$ws.Cells.Item(141, 8).Value = 76701.664
$ws.Cells.Item(141, 10).Value = 76701.664
$ws.Cells.Item(141, 12).Value = 76701.664
$ws.Cells.Item(141, 14).Value = -87061.664

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 986.8182
$ws.Cells.Item(16, 9).Value = 930.25
$ws.Cells.Item(16, 10).Value = 1137.6666
$ws.Cells.Item(16, 11).Value = 930.25
$ws.Cells.Item(16, 12).Value = 1137.6666
$ws.Cells.Item(16, 13).Value = -643.25
$ws.Cells.Item(16, 14).Value = -1711.6666
$ws.Cells.Item(113, 8).Value = 986.8182
$ws.Cells.Item(113, 9).Value = 930.25
$ws.Cells.Item(113, 10).Value = 1137.6666
$ws.Cells.Item(113, 11).Value = 930.25
$ws.Cells.Item(113, 12).Value = 1137.6666
$ws.Cells.Item(113, 13).Value = 1239.75
$ws.Cells.Item(113, 14).Value = -5477.6666
$ws.Cells.Item(132, 8).Value = 1424.1552
$ws.Cells.Item(132, 9).Value = 1091.8536
$ws.Cells.Item(132, 10).Value = 2225.5881
$ws.Cells.Item(132, 11).Value = 3275.5608
$ws.Cells.Item(132, 12).Value = 6676.7643
$ws.Cells.Item(132, 13).Value = -745.5607999999997
$ws.Cells.Item(132, 14).Value = -11736.7643

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 4200
$ws.Cells.Item(41, 10).Value = 9250
$ws.Cells.Item(41, 12).Value = 9250
$ws.Cells.Item(41, 14).Value = -9960
$ws.Cells.Item(70, 8).Value = 6849393
$ws.Cells.Item(70, 9).Value = 2721698
$ws.Cells.Item(70, 10).Value = 20411820
$ws.Cells.Item(70, 11).Value = 2721698
$ws.Cells.Item(70, 12).Value = 20411820
$ws.Cells.Item(70, 13).Value = -2721428
$ws.Cells.Item(70, 14).Value = -20412360
$ws.Cells.Item(73, 8).Value = 6849393
$ws.Cells.Item(73, 9).Value = 2721698
$ws.Cells.Item(73, 10).Value = 20411820
$ws.Cells.Item(73, 11).Value = 2721698
$ws.Cells.Item(73, 12).Value = 20411820
$ws.Cells.Item(73, 13).Value = -2720762
$ws.Cells.Item(73, 14).Value = -20413692
$ws.Cells.Item(113, 8).Value = 23860.916
$ws.Cells.Item(113, 9).Value = 1069.6666
$ws.Cells.Item(113, 10).Value = 46652.168
$ws.Cells.Item(113, 11).Value = 1069.6666
$ws.Cells.Item(113, 12).Value = 46652.168
$ws.Cells.Item(113, 13).Value = 1100.3334
$ws.Cells.Item(113, 14).Value = -50992.168

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 50003600
$ws.Cells.Item(46, 9).Value = 3500
$ws.Cells.Item(46, 10).Value = 71432216
$ws.Cells.Item(46, 11).Value = 3500
$ws.Cells.Item(46, 12).Value = 71432216
$ws.Cells.Item(46, 13).Value = -3312
$ws.Cells.Item(46, 14).Value = -71432592
$ws.Cells.Item(61, 8).Value = 4355.3335
$ws.Cells.Item(61, 9).Value = 3166.2
$ws.Cells.Item(61, 11).Value = 3166.2
$ws.Cells.Item(61, 13).Value = -2964.2
$ws.Cells.Item(68, 8).Value = 1880.9445
$ws.Cells.Item(68, 9).Value = 1196.7273
$ws.Cells.Item(68, 10).Value = 2956.1428
$ws.Cells.Item(68, 11).Value = 1196.7273
$ws.Cells.Item(68, 12).Value = 2956.1428
$ws.Cells.Item(68, 13).Value = -447.7273
$ws.Cells.Item(68, 14).Value = -4454.1428
$ws.Cells.Item(71, 8).Value = 1880.9445
$ws.Cells.Item(71, 9).Value = 1196.7273
$ws.Cells.Item(71, 10).Value = 2956.1428
$ws.Cells.Item(71, 11).Value = 5983.636500000001
$ws.Cells.Item(71, 12).Value = 14780.714
$ws.Cells.Item(71, 13).Value = -2239.636500000001
$ws.Cells.Item(71, 14).Value = -22268.714
$ws.Cells.Item(113, 8).Value = 4355.3335
$ws.Cells.Item(113, 9).Value = 3166.2
$ws.Cells.Item(113, 11).Value = 3166.2
$ws.Cells.Item(113, 13).Value = -996.1999999999998

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 8209.808000000001
$ws.Cells.Item(107, 9).Value = 13609.533
$ws.Cells.Item(107, 11).Value = 40828.599
$ws.Cells.Item(107, 13).Value = -38908.599
$ws.Cells.Item(113, 8).Value = 121
$ws.Cells.Item(113, 9).Value = 95
$ws.Cells.Item(113, 10).Value = 147
$ws.Cells.Item(113, 11).Value = 285
$ws.Cells.Item(113, 12).Value = 441
$ws.Cells.Item(113, 13).Value = 1885
$ws.Cells.Item(113, 14).Value = -4781
$ws.Cells.Item(122, 8).Value = 1535.381
$ws.Cells.Item(122, 9).Value = 1261.9166
$ws.Cells.Item(122, 10).Value = 1900
$ws.Cells.Item(122, 11).Value = 3785.7498
$ws.Cells.Item(122, 12).Value = 5700
$ws.Cells.Item(122, 13).Value = -1335.7498
$ws.Cells.Item(122, 14).Value = -10600
$ws.Cells.Item(126, 8).Value = 10000798
$ws.Cells.Item(126, 9).Value = 13889520
$ws.Cells.Item(126, 10).Value = 1227
$ws.Cells.Item(126, 11).Value = 41668560
$ws.Cells.Item(126, 12).Value = 3681
$ws.Cells.Item(126, 13).Value = -41666090
$ws.Cells.Item(126, 14).Value = -8621
$ws.Cells.Item(132, 8).Value = 2076405.2
$ws.Cells.Item(132, 9).Value = 7504.143
$ws.Cells.Item(132, 11).Value = 22512.429
$ws.Cells.Item(132, 13).Value = -19982.429
